$d = $word.ActiveDocument

# Step 1: collapse the sentence's text content to its final wording while it
# is still a single run - this is the simplest way to get the right text in
# place, it will happen to re-merge surrounding runs, which is fine since
# none of them are touched here.
$para = $d.Paragraphs.Item(1)
$rng = $para.Range.Duplicate
$rng.Find.ClearFormatting()
$rng.Find.Execute("de référence à une variable", $false, $false, $false, $false, $false, $true, 1, $false, "conditionnelles", 2)

# Step 2: split "Template de test pour les balises conditionnelles : " into
# three separate runs, matching the authored edit:
#   "Template de test pour les balises " | "conditionnelles" | " : "
# Do this LAST (after all text-content edits) via a bookmark add/delete
# trick, which breaks run boundaries cleanly (no stray rPr) without
# triggering the engine's run-merge normalisation that follows any direct
# text-content mutation.
$para2 = $d.Paragraphs.Item(1)
$find2 = $para2.Range.Duplicate
$find2.Find.ClearFormatting()
$found = $find2.Find.Execute("conditionnelles")
if (-not $found) {
    throw "could not locate replaced phrase"
}

$splitStart = $find2.Start
$splitEnd = $find2.End

$b1 = $d.Range($splitStart, $splitStart)
$d.Bookmarks.Add("zzSplitA", $b1)
$b2 = $d.Range($splitEnd, $splitEnd)
$d.Bookmarks.Add("zzSplitB", $b2)

$d.Bookmarks("zzSplitA").Delete()
$d.Bookmarks("zzSplitB").Delete()
